# option.xlsx — add the missing "Interact" event row to the Sheet2 table.
#
# Sheet2 (xl/worksheets/sheet1.xml) is a simple Id / Name / EventName lookup
# table. This bugfix appends row 8 for the "交互" (Interact) event, which was
# missing: Id=1011, Name="交互", EventName="OnInteractEvent".

# A few VBA/Excel constants (no interop type library loaded in this host,
# so use their well-known numeric values directly).
$xlPasteFormats = -4122
$xlLeft = -4131
$xlCenter = -4108

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Id) -----------------------------------------------------
# Reuse the existing "Id" column formatting (same style as A3:A7) by
# copying format from the row above, then write the new id.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = 1011

# --- Column B (Name) ----------------------------------------------------
# Write the Chinese label first so it lands on the lower shared-string
# index (matching insertion order left-to-right across the new row).
$ws.Range("B8").Value = "交互"

# --- Column C (EventName) ------------------------------------------------
# Reuse the existing "EventName" column formatting (same style as C3:C7)
# by copying format from the row above, then write the new event name.
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C8").Value = "OnInteractEvent"

# Give B8 left/center aligned text matching the rest of the table (font
# formatting applied explicitly, same as the other Name cells visually).
$ws.Range("B8").Font.Name = $ws.Range("C4").Font.Name
$ws.Range("B8").Font.Size = $ws.Range("C4").Font.Size
$ws.Range("B8").HorizontalAlignment = $xlLeft
$ws.Range("B8").VerticalAlignment = $xlCenter

# Move the window / selection to the newly added cell, like the author did.
$ws.Range("C8").Select() | Out-Null
try {
    $win = $excel.ActiveWindow
    $win.Left = 0
    $win.Top = 2865
} catch {
}
